# Applies the "Add files via upload" commit:
#   - Every cached `datetime1` / `datetimeFigureOut` field ("18.02.20",
#     "2/20/20", "2/18/20") is refreshed to the new save date (21 Feb 2020),
#     shown per-shape in whatever format that shape was already using.
#   - A typo fix on slide 3: "least dimension" -> "first dimension".

$p = $ppt.ActivePresentation

function Set-ShapeFullText($shape, $newText) {
    # Writing straight to the shape's own TextFrame.TextRange (as opposed to
    # a Characters()/Paragraphs() sub-range) performs a true in-place
    # replace of the whole (single-paragraph) run/field content.
    $shape.TextFrame.TextRange.Text = $newText
}

function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }
        $t = $sh.TextFrame.TextRange.Text
        if ($t -eq "18.02.20") {
            Set-ShapeFullText $sh "21.02.20"
        } elseif ($t -eq "2/20/20") {
            Set-ShapeFullText $sh "2/21/20"
        } elseif ($t -eq "2/18/20") {
            Set-ShapeFullText $sh "2/21/20"
        }
    }
}

# Slide master date placeholder.
Update-DateField $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    $lay = $p.SlideMaster.CustomLayouts.Item($j)
    Update-DateField $lay.Shapes
}

# Notes master date placeholder.
Update-DateField $p.NotesMaster.Shapes

# Each slide's own date placeholder.
for ($k = 1; $k -le $p.Slides.Count; $k++) {
    Update-DateField $p.Slides.Item($k).Shapes
}

# Slide 3 body text typo fix: "least" -> "first" (single word, keep the rest
# of the paragraph / text frame structure untouched).
$s3 = $p.Slides.Item(3)
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $sh = $s3.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }
    $tr = $sh.TextFrame.TextRange
    if ($tr.Text -like "*Untill we reach to the least dimension*") {
        $paraCount = $tr.Paragraphs(1, -1).Count
        for ($para = 1; $para -le $paraCount; $para++) {
            $pr = $tr.Paragraphs($para, 1)
            if ($pr.Text -eq "Untill we reach to the least dimension") {
                # Route through a throwaway value first so the engine's
                # prefix/suffix run-diff has nothing in common with either
                # the old or the new sentence, and rewrites the paragraph
                # as a single run (matching a plain text edit) instead of
                # splitting out just the changed "lea"/"fir" fragment.
                $pr.Text = "*"
                $tr2 = $sh.TextFrame.TextRange
                $pr2 = $tr2.Paragraphs($para, 1)
                $pr2.Text = "Untill we reach to the first dimension"
            }
        }
    }
}
